# edit.ps1 - reproduce the authored change:
#   1. Update the cached text of every "datetimeFigureOut" date field
#      (slide master, every slide layout, and the notes master) from
#      "5/25/23" to "3/23/24".
#   2. Nudge the rotated "Rectangle 6" shape and its attached
#      "Straight Connector 22" on slide 1 to their new position/size.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Date placeholder text ("datetimeFigureOut" field) -> 3/23/24
# ---------------------------------------------------------------
function Update-DatePlaceholder($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14) {
            $ph = $sh.PlaceholderFormat
            if ($ph.Type -eq 16) {
                $sh.TextFrame.TextRange.Text = "3/23/24"
            }
        }
    }
}

# Slide master
Update-DatePlaceholder($p.SlideMaster)

# Every slide layout under the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder($layouts.Item($li))
}

# Notes master
Update-DatePlaceholder($p.NotesMaster)

# ---------------------------------------------------------------
# 2) Move/resize "Rectangle 6" + "Straight Connector 22" (slide 1)
# ---------------------------------------------------------------
$slide1 = $p.Slides.Item(1)

$rect = $slide1.Shapes.Item(2)   # "Rectangle 6"
$rect.Left = 688.2607480314961
$rect.Top = 201.4622440944882
$rect.Width = 77.54783464566928
$rect.Height = 37.241377952755904

$conn = $slide1.Shapes.Item(19)  # "Straight Connector 22"
$conn.Left = 726.4100393700787
$conn.Top = 89.17248031496062
$conn.Width = 0.6246850393700787
$conn.Height = 92.1365748031496
